# Auto update: 2025-12-05 12:21:03
# Refresh the gold-hedge watchlist: row 2 (GLD/StreetTRACKS) and row 3
# (GC=F/Gold Feb 26) swap places and pick up refreshed market data, and
# the MACRO_SCORE column (N) is recomputed for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> StreetTRACKS Gold Shares / GLD
$ws.Range("B2").Value = "StreetTRACKS Gold Shares"
$ws.Range("C2").Value = "GLD"
$ws.Range("D2").Value = 387.13
$ws.Range("E2").Value = 56.3
$ws.Range("F2").Value = 1.05
$ws.Range("H2").Value = 73
$ws.Range("J2").Value = 96
$ws.Range("N2").Value = 54.84087454262382

# Row 3 -> Gold Feb 26 / GC=F
$ws.Range("B3").Value = "Gold Feb 26"
$ws.Range("C3").Value = "GC=F"
$ws.Range("D3").Value = 4232.2
$ws.Range("E3").Value = 55.1
$ws.Range("F3").Value = 1.61
$ws.Range("H3").Value = 63
$ws.Range("I3").Value = 80
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 66.5
$ws.Range("N3").Value = 54.84087454262382

# Row 4 (Newmont Corporation / NEM) keeps its values, only the
# recomputed MACRO_SCORE changes.
$ws.Range("N4").Value = 54.84087454262382
